# Adds a new "CHARACTERS" section with a two-item bulleted list to the end
# of the design doc, right after the paragraph that ends with
# "...possibly opening up new rp options through them."

$d = $word.ActiveDocument

# Locate the anchor paragraph (the last paragraph in the body).
$found = $d.Content.Find.Execute(
    "possibly opening up new rp options through them.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchor = $d.Paragraphs.Last.Range
$anchor.Collapse(0)

# Blank separator paragraph.
$anchor.InsertParagraphAfter() | Out-Null
$anchor.Collapse(0)

# "CHARACTERS" heading paragraph.
$anchor.InsertParagraphAfter() | Out-Null
$anchor.Collapse(0)
$headingPara = $d.Paragraphs.Last
$headingPara.Range.InsertBefore("CHARACTERS") | Out-Null

# First bullet: Miskovine.
$anchor2 = $d.Paragraphs.Last.Range
$anchor2.Collapse(0)
$anchor2.InsertParagraphAfter() | Out-Null
$anchor2.Collapse(0)
$bullet1 = $d.Paragraphs.Last
$bullet1.Range.InsertBefore("Miskovine, boss of the Ashfoot gang") | Out-Null
$bullet1.Range.ListFormat.ApplyBulletDefault() | Out-Null

# Second bullet: Kellach (continues the same bulleted list).
$anchor3 = $d.Paragraphs.Last.Range
$anchor3.Collapse(0)
$anchor3.InsertParagraphAfter() | Out-Null
$anchor3.Collapse(0)
$bullet2 = $d.Paragraphs.Last
$bullet2.Range.InsertBefore("Kellach, traitor and double agent to the Ashfoot gang") | Out-Null
$bullet2.Range.ListFormat.ApplyBulletDefault() | Out-Null
